# Auto-generated Excel COM-interop script
# Applies the "cryptos" price/volume update described by the commit diff:
# refreshed D (Price) / E (Volume 1h) figures for most rows, plus a pair
# of row swaps (ShibaInu<->Avalanche at rows 12/13, EnergySwap<->Algorand
# at rows 43/44) whose B/C/D/E cells were all rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $val) {
    # Source cells are plain text (inline strings). Values that look like a
    # bare integer/decimal number (e.g. "45.53", "1.00") would otherwise be
    # auto-coerced to a Number by Excel on assignment, so prefix those with
    # an apostrophe - the standard Excel "force text" marker - which keeps
    # the stored value as clean text (just flips on the quote-prefix flag)
    # without altering values that are already unambiguously text (coin
    # names, URLs, percent strings, multi-dot price strings like
    # "97.669.71").
    if ($val -match '^\d+(\.\d+)?$') {
        $ws.Range($cellRef).Value = "'" + $val
    } else {
        $ws.Range($cellRef).Value = $val
    }
}

Set-CellText 'D2' '97.669.71'
Set-CellText 'E2' '  +1.33%  '

Set-CellText 'D3' '3.719.75'
Set-CellText 'E3' '  +0.11%  '

Set-CellText 'D4' '0.999'
Set-CellText 'E4' '  -0.17%  '

Set-CellText 'D5' '2.18'
Set-CellText 'E5' '  +13.08%  '

Set-CellText 'D6' '238.39'
Set-CellText 'E6' '  +0.06%  '

Set-CellText 'D7' '657.83'
Set-CellText 'E7' '  +0.55%  '

Set-CellText 'D8' '0.445'
Set-CellText 'E8' '  +5.47%  '

Set-CellText 'D9' '1.14'
Set-CellText 'E9' '  +6.05%  '

Set-CellText 'E10' '  -0.11%  '

Set-CellText 'D11' '3.719.38'
Set-CellText 'E11' '  +0.17%  '

Set-CellText 'B12' 'Avalanche'
Set-CellText 'C12' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-CellText 'D12' '45.53'
Set-CellText 'E12' '  +1.31%  '

Set-CellText 'B13' 'ShibaInu'
Set-CellText 'C13' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-CellText 'D13' '0.0000312'
Set-CellText 'E13' '  +16.96%  '

Set-CellText 'E14' '  +0.93%  '

Set-CellText 'D15' '6.86'
Set-CellText 'E15' '  +0.21%  '

Set-CellText 'D16' '4.413.75'
Set-CellText 'E16' '  +0.08%  '

Set-CellText 'D17' '97.401.07'
Set-CellText 'E17' '  +1.20%  '

Set-CellText 'E18' '  -1.12%  '

Set-CellText 'D19' '3.729.54'
Set-CellText 'E19' '  +0.10%  '

Set-CellText 'D20' '13.15'
Set-CellText 'E20' '  +2.97%  '

Set-CellText 'D21' '18.93'
Set-CellText 'E21' '  -0.84%  '

Set-CellText 'D22' '0.540'
Set-CellText 'E22' '  +2.91%  '

Set-CellText 'D23' '533.05'
Set-CellText 'E23' '  +1.83%  '

Set-CellText 'E24' '  -0.02%  '

Set-CellText 'D25' '0.0000224'
Set-CellText 'E25' '  +10.75%  '

Set-CellText 'D26' '119.44'
Set-CellText 'E26' '  +16.69%  '

Set-CellText 'E27' '  -1.50%  '

Set-CellText 'D28' '0.217'
Set-CellText 'E28' '  +29.04%  '

Set-CellText 'E29' '  +0.38%  '

Set-CellText 'E30' '  +3.39%  '

Set-CellText 'D31' '3.06'
Set-CellText 'E31' '  -0.21%  '

Set-CellText 'D32' '1.00'
Set-CellText 'E32' '  +0.07%  '

Set-CellText 'D33' '0.192'
Set-CellText 'E33' '  +3.35%  '

Set-CellText 'E34' '  -2.34%  '

Set-CellText 'D35' '33.22'
Set-CellText 'E35' '  +1.36%  '

Set-CellText 'D36' '0.997'
Set-CellText 'E36' '  -0.76%  '

Set-CellText 'E37' '  +0.74%  '

Set-CellText 'D38' '640.95'
Set-CellText 'E38' '  -3.96%  '

Set-CellText 'E39' '  -0.84%  '

Set-CellText 'E41' '  +5.29%  '

Set-CellText 'D42' '6.91'
Set-CellText 'E42' '  -2.76%  '

Set-CellText 'B43' 'Algorand'
Set-CellText 'C43' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-CellText 'D43' '0.497'
Set-CellText 'E43' '  +12.88%  '

Set-CellText 'B44' 'EnergySwap'
Set-CellText 'C44' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText 'D44' '41.12'
Set-CellText 'E44' '  +1.49%  '

Set-CellText 'D45' '2.03'
Set-CellText 'E45' '  +3.05%  '

Set-CellText 'D46' '0.971'
Set-CellText 'E46' '  -0.79%  '

Set-CellText 'D47' '0.0463'
Set-CellText 'E47' '  +0.89%  '

Set-CellText 'D48' '2.41'
Set-CellText 'E48' '  +3.82%  '

Set-CellText 'D49' '8.96'
Set-CellText 'E49' '  +4.30%  '

Set-CellText 'E50' '  +0.24%  '

Set-CellText 'E51' '  +6.14%  '
